$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Attendance table: Eda (row 7) now answers "ANO" instead of "NE" ---
# Copy Vojta's (row 6) "ANO" formatting onto Eda's answer cell, then set the value.
$ws.Range("B6").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B7").Value = "ANO"

# --- "Navrhy dnu a casu" table: highlight the time slot for 19.10. and add new rows ---
# D15 (the "navecer" time cell for 19.10.) gets highlighted the same as C15.
$ws.Range("C15").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New row 16: another slot on 19.10., "cely den" (whole day)
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null
$ws.Range("A16").Value = "19.10."

$ws.Range("B15").Copy() | Out-Null
$ws.Range("B16").PasteSpecial(-4122) | Out-Null
$ws.Range("B16").Value = "cel" + [char]0x00FD + " den"

$ws.Range("D15").Copy() | Out-Null
$ws.Range("D16").PasteSpecial(-4122) | Out-Null

# New row 17: 18.10., "cely den"
$ws.Range("A17").Value = "18.10."
$ws.Range("B17").Value = "cel" + [char]0x00FD + " den"
$ws.Range("D15").Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4122) | Out-Null

# New row 18: 17.10., "cely den"
$ws.Range("A18").Value = "17.10."
$ws.Range("B18").Value = "cel" + [char]0x00FD + " den"
$ws.Range("D15").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4122) | Out-Null

# --- Move the active selection to D18 ---
$ws.Range("D18").Select() | Out-Null

Write-Output "done"
